$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "26.596.06"
Set-TextValue "E2" "  +1.15%  "

Set-TextValue "D3" "1.630.55"
Set-TextValue "E3" "  +1.37%  "

Set-TextValue "D5" "212.58"
Set-TextValue "E5" "  -0.17%  "

Set-TextValue "B6" "XRP"
Set-TextValue "C6" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D6" "0.493"
Set-TextValue "E6" "  +1.31%  "

Set-TextValue "B7" "USDC"
Set-TextValue "C7" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  +0.05%  "

Set-TextValue "D8" "0.252"
Set-TextValue "E8" "  +0.77%  "

Set-TextValue "E9" "  +1.47%  "

Set-TextValue "E10" "  +2.65%  "

Set-TextValue "D11" "0.0839"
Set-TextValue "E11" "  +3.18%  "

Set-TextValue "D12" "1.858.55"
Set-TextValue "E12" "  +1.41%  "

Set-TextValue "D13" "1.636.39"
Set-TextValue "E13" "  +1.73%  "

Set-TextValue "D14" "4.08"
Set-TextValue "E14" "  +1.32%  "

Set-TextValue "E15" "  +2.02%  "

Set-TextValue "D16" "26.597.98"

Set-TextValue "E17" "  +1.21%  "

Set-TextValue "E18" "  +1.38%  "

Set-TextValue "D19" "209.33"

Set-TextValue "E20" "  -0.01%  "

Set-TextValue "E21" "  +0.64%  "

Set-TextValue "E22" "  +0.78%  "

Set-TextValue "D23" "6.18"
Set-TextValue "E23" "  +2.74%  "

Set-TextValue "D24" "1.93"
Set-TextValue "E24" "  +2.58%  "

Set-TextValue "D25" "146.84"
Set-TextValue "E25" "  +2.47%  "

Set-TextValue "E26" "  +0.02%  "

Set-TextValue "E27" "  -0.69%  "

Set-TextValue "E28" "  +4.13%  "

Set-TextValue "D29" "15.33"
Set-TextValue "E29" "  +0.60%  "

Set-TextValue "D30" "0.0523"
Set-TextValue "E30" "  +5.12%  "

Set-TextValue "E31" "  -0.31%  "

Set-TextValue "D32" "3.24"
Set-TextValue "E32" "  +1.29%  "

Set-TextValue "E33" "  -0.03%  "

Set-TextValue "E34" "  +0.69%  "

Set-TextValue "E35" "  +0.04%  "

Set-TextValue "D36" "1.164.79"
Set-TextValue "E36" "  +0.26%  "

Set-TextValue "E37" "  +0.47%  "

Set-TextValue "D38" "0.807"
Set-TextValue "E38" "  +2.31%  "

Set-TextValue "E39" "  +0.04%  "

Set-TextValue "D40" "0.503"
Set-TextValue "E40" "  +1.42%  "

Set-TextValue "E41" "  -0.27%  "

Set-TextValue "D42" "0.793"
Set-TextValue "E42" "  +1.11%  "

Set-TextValue "E43" "  +0.08%  "

Set-TextValue "D44" "1.771.44"
Set-TextValue "E44" "  +1.59%  "

Set-TextValue "D45" "91.97"
Set-TextValue "E45" "  -0.08%  "

Set-TextValue "E46" "  +0.59%  "

Set-TextValue "D47" "54.58"
Set-TextValue "E47" "  +0.98%  "

Set-TextValue "E48" "  -2.67%  "

Set-TextValue "E49" "  +0.68%  "

Set-TextValue "E50" "  +4.48%  "

Set-TextValue "E51" "  +0.50%  "

